$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D. Excel shifts the existing
# D:K quarterly data right to F:M, and the new D:E columns are ready
# to receive the two newest quarters of data.
$ws.Columns("D:E").Insert(-4161, 1)

# The freshly inserted D:E columns come back with the default/general
# format instead of the per-row format (date row vs plain-number rows)
# used throughout the rest of the table. Clone that formatting from the
# (now shifted) original D:E columns, which live at F:G after the insert.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the column widths used elsewhere in the table.
$ws.Columns("D").ColumnWidth = $ws.Columns("G").ColumnWidth
$ws.Columns("E").ColumnWidth = $ws.Columns("I").ColumnWidth

# Populate the two new quarter columns (D = newest quarter, E = the
# quarter before it) with the values from this update.
$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404
$ws.Range("D8").Value = 42000
$ws.Range("E8").Value = 14700
$ws.Range("D9").Value = 40000
$ws.Range("E9").Value = 18200
$ws.Range("D10").Value = 2000
$ws.Range("E10").Value = -3500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = 1600
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 45000
$ws.Range("E17").Value = 24300
$ws.Range("D18").Value = -3000
$ws.Range("E18").Value = -9600
$ws.Range("D20").Value = -3400
$ws.Range("E20").Value = 4700
$ws.Range("D21").Value = -4300
$ws.Range("E21").Value = -2800
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("D23").Value = -6400
$ws.Range("E23").Value = -4900
$ws.Range("D24").Value = -1800
$ws.Range("E24").Value = -1300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -4700
$ws.Range("E26").Value = -3500
$ws.Range("D27").Value = -4800
$ws.Range("E27").Value = -3700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = 300
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 3400
$ws.Range("E32").Value = -4700
$ws.Range("D33").Value = -4800
$ws.Range("E33").Value = -3400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -4800
$ws.Range("E35").Value = -3400
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404
$ws.Range("D41").Value = 1300
$ws.Range("E41").Value = 600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 22500
$ws.Range("E43").Value = 17300
$ws.Range("D44").Value = 1300
$ws.Range("E44").Value = 1200
$ws.Range("D45").Value = 11600
$ws.Range("E45").Value = 11900
$ws.Range("D46").Value = 36700
$ws.Range("E46").Value = 31000
$ws.Range("D47").Value = 82200
$ws.Range("E47").Value = 48500
$ws.Range("D48").Value = 245600
$ws.Range("E48").Value = 332800
$ws.Range("D49").Value = 7800
$ws.Range("E49").Value = 7700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1100
$ws.Range("E52").Value = 1300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 373500
$ws.Range("E54").Value = 421300
$ws.Range("D57").Value = 8000
$ws.Range("E57").Value = 6100
$ws.Range("D58").Value = 2900
$ws.Range("E58").Value = 3100
$ws.Range("D59").Value = 15800
$ws.Range("E59").Value = 17800
$ws.Range("D60").Value = 26700
$ws.Range("E60").Value = 27100
$ws.Range("D61").Value = 94000
$ws.Range("E61").Value = 77000
$ws.Range("D62").Value = 27500
$ws.Range("E62").Value = 87300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 148800
$ws.Range("E66").Value = 192000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 10800
$ws.Range("E70").Value = 10800
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 58400
$ws.Range("E72").Value = 50400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 213900
$ws.Range("E76").Value = 218600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404
$ws.Range("D81").Value = -4800
$ws.Range("E81").Value = -3400
$ws.Range("D83").Value = 2100
$ws.Range("E83").Value = 2100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -4500
$ws.Range("E89").Value = -4800
$ws.Range("D91").Value = -400
$ws.Range("E91").Value = -4100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -9800
$ws.Range("E94").Value = 0
$ws.Range("D96").Value = -1300
$ws.Range("E96").Value = -1100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 14900
$ws.Range("E100").Value = 4900
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 700
$ws.Range("E102").Value = 100

# Rows that only held a section label (in column A or B) had no cells
# at all in D:K originally, so they should not end up with empty,
# merely-formatted D:E cells either - clear those back out.
$ws.Range("D5:E6").Clear()
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()
